# EI Variable Installments T1 scenarios
# - Recompute Summary "Fees" column (row 5) and drop the now-empty trailing row.
# - Recompute Repayment schedule Fees/Due/Outstanding for rows 3-5.
# - Active tab moves from NewLoanInput to Transactions.

$wb = $excel.ActiveWorkbook

# --- Summary sheet -------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Summary")
$ws2.Range("A5").Value = 16.79
$ws2.Range("E5").Value = 16.79
$ws2.Range("F5").Value = 16.79
# Row 6 (all zeros) is no longer part of the data set - remove it.
$ws2.Rows("6:16").Select()
$ws2.Rows("6").Delete()

# --- Repayment schedule sheet ---------------------------------------------
$ws3 = $wb.Worksheets.Item("Repayment schedule")
$ws3.Range("J3").Value = 8.3699999999999992
$ws3.Range("K3").Value = 896.09
$ws3.Range("P3").Value = 896.09
$ws3.Range("J4").Value = 8.42
$ws3.Range("K4").Value = 896.14
$ws3.Range("P4").Value = 896.14
$ws3.Range("J5").Value = 0
$ws3.Range("K5").Value = 887.72
$ws3.Range("P5").Value = 887.72
$ws3.Range("F17").Select()

# --- Transactions becomes the active tab/sheet -----------------------------
$ws4 = $wb.Worksheets.Item("Transactions")
$ws4.Activate()
$ws4.Range("I9").Select()
